$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reference element/parameter block (row 2-4): atomic symbol,
# xi and di values for the new element (Pt -> Si).
$ws.Range("N2").Value = "Si"
$ws.Range("N3").Value = 4.295
$ws.Range("N4").Value = 0.402

# Add the LENNARD-JONES RCUT value (10) for each of the seven parameter
# blocks that previously had no RCUT entry.
$ws.Range("Q5").Value = 10
$ws.Range("Q12").Value = 10
$ws.Range("Q19").Value = 10
$ws.Range("Q26").Value = 10
$ws.Range("Q33").Value = 10
$ws.Range("Q40").Value = 10
$ws.Range("Q47").Value = 10

# Move the selection / view back to the top of the sheet.
$ws.Range("M7").Select()
